# WAY Corrector - تحديث الملف: 21/01/2026, 16:27:11
# Append two new rows (4 and 5) to the sheet, mirroring the existing
# rows 1-3: column A holds a text value, column B holds an (ignored,
# number-stored-as-text) empty text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the new cells keep their values as text (matches the
# "numberStoredAsText" pattern already used by A1:B3).
$ws.Range("A4:B5").NumberFormat = "@"

$ws.Range("A4").Value = "679056"
$ws.Range("B4").Value = "'"

$ws.Range("A5").Value = "ztuome"
$ws.Range("B5").Value = "'"
